$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '35.288.40'
$ws.Range("E2").Value = '  +0.44%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.878.94'
$ws.Range("E3").Value = '  -1.41%  '

$ws.Range("E4").Value = '  -0.49%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '245.31'
$ws.Range("E5").Value = '  -3.32%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.683'
$ws.Range("E6").Value = '  -2.05%  '

$ws.Range("E7").Value = '  -0.58%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '43.14'
$ws.Range("E8").Value = '  +3.98%  '

$ws.Range("E9").Value = '  -1.52%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '53.38'
$ws.Range("E10").Value = '  +1.62%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0738'
$ws.Range("E11").Value = '  -1.94%  '

$ws.Range("E12").Value = '  -0.84%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '13.44'
$ws.Range("E13").Value = '  +1.83%  '

$ws.Range("E14").Value = '  -1.31%  '

$ws.Range("E16").Value = '  -2.62%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.860.52'
$ws.Range("E17").Value = '  -2.38%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '35.300.08'

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '72.72'
$ws.Range("E19").Value = '  -1.43%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0818'
$ws.Range("E20").Value = '  -2.55%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '243.64'
$ws.Range("E21").Value = '  +0.25%  '

$ws.Range("E22").Value = '  -1.81%  '

$ws.Range("E23").Value = '  -2.12%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.62'
$ws.Range("E24").Value = '  +7.57%  '

$ws.Range("E25").Value = '  -0.55%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.15'
$ws.Range("E26").Value = '  -6.95%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '165.16'
$ws.Range("E27").Value = '  -1.63%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.52'
$ws.Range("E28").Value = '  -0.84%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '18.23'
$ws.Range("E29").Value = '  -1.76%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.126'
$ws.Range("E30").Value = '  -2.55%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.128.44'
$ws.Range("E31").Value = '  +0.00%  '

$ws.Range("E32").Value = '  +7.41%  '

$ws.Range("B33").Value = 'WEMIXToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.99'
$ws.Range("E33").Value = '  -1.11%  '

$ws.Range("B34").Value = 'Filecoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.25'
$ws.Range("E34").Value = '  -2.09%  '

$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0585'
$ws.Range("E35").Value = '  -3.04%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.12'
$ws.Range("E36").Value = '  -2.18%  '

$ws.Range("E37").Value = '  -0.55%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.839'
$ws.Range("E38").Value = '  -1.54%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0741'
$ws.Range("E39").Value = '  +13.83%  '

$ws.Range("E40").Value = '  -4.08%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '17.61'
$ws.Range("E41").Value = '  +2.00%  '

$ws.Range("E42").Value = '  +0.01%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '96.13'
$ws.Range("E43").Value = '  -7.01%  '

$ws.Range("E44").Value = '  -2.93%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.302.54'
$ws.Range("E45").Value = '  -0.11%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.35'
$ws.Range("E46").Value = '  -2.81%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0795'
$ws.Range("E47").Value = '  +5.89%  '

$ws.Range("E48").Value = '  -1.72%  '

$ws.Range("E49").Value = '  -1.29%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '12.12'
$ws.Range("E50").Value = '  -4.46%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.21'
$ws.Range("E51").Value = '  -5.80%  '
